$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

# --- New tool column added to the right of the table: "Trisotech BPMN Visio
#     Modeler" version 4.2, contributed as BI:BK (a 3-column group matching the
#     layout used by the other tools-under-test columns) ---
$ws.Range("BI3").Value = "Trisotech BPMN Visio Modeler"
$ws.Range("BI4").Value = "4.2"
$ws.Range("BI5").Value = 41716

# --- Test results ("OK") reported for the new tool across the existing test
#     rows. The BI (and, for a couple of rows, BJ) cells still carried an
#     incomplete border left over from the template, so once they get a
#     value their format is picked up from the neighbouring BK cell in the
#     same row, which already has the finished look. ---
$ws.Range("BI6").Value = "OK"
$ws.Range("BK6").Copy()
$ws.Range("BI6").PasteSpecial($xlPasteFormats)
$ws.Range("BJ6").Value = "OK"
$ws.Range("BK6").Value = "OK"

$ws.Range("BI7").Value = "OK"
$ws.Range("BK7").Copy()
$ws.Range("BI7").PasteSpecial($xlPasteFormats)
$ws.Range("BJ7").Value = "OK"
$ws.Range("BK7").Value = "OK"

$ws.Range("BI8").Value = "OK"
$ws.Range("BK8").Copy()
$ws.Range("BI8").PasteSpecial($xlPasteFormats)
$ws.Range("BJ8").Value = "OK"
$ws.Range("BK8").Value = "OK"

$ws.Range("BI9").Value = "OK"
$ws.Range("BK9").Copy()
$ws.Range("BI9").PasteSpecial($xlPasteFormats)
$ws.Range("BJ9").PasteSpecial($xlPasteFormats)
$ws.Range("BJ9").Value = "OK"
$ws.Range("BK9").Value = "OK"

$ws.Range("AT10").Value = "OK"
$ws.Range("AU10").Value = "OK"
$ws.Range("AV10").Value = "OK"

$ws.Range("BI10").Value = "OK"
$ws.Range("BK10").Copy()
$ws.Range("BI10").PasteSpecial($xlPasteFormats)
$ws.Range("BJ10").PasteSpecial($xlPasteFormats)
$ws.Range("BJ10").Value = "OK"
$ws.Range("BK10").Value = "OK"

$ws.Range("BI11").Value = "OK"
$ws.Range("BJ11").Value = "OK"
$ws.Range("BK11").Value = "OK"

$ws.Range("BI12").Value = "OK"
$ws.Range("BJ12").Value = "OK"
$ws.Range("BK12").Value = "OK"
